# Change table to nanoseconds per day
#
# The second (25K) table in Sheet1 only went up through row 15 (Wall Time
# in hours). Extend it with a second "ns/day" style block in column I that
# converts each of those wall-time values (I11:I15) to the same units as
# column I's first table (minutes -> the *60 conversion used elsewhere in
# the sheet), placed two rows below the existing block (I17:I21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I17 = I11 * 60 (single formula, own <f> element)
$ws.Range("I17").Formula = "=I11*60"

# I18:I21 = I12*60 .. I15*60, entered as one relative formula across the
# range so Excel stores it as a shared formula (t="shared") the way the
# original table's H column does.
$ws.Range("I18:I21").Formula = "=I12*60"

# Update the view: the selection moves to the first (25K) Wall Time column
# block, and the sheet scrolls back up so A1 is visible again (no more
# topLeftCell scroll-offset).
$ws.Range("H11:H15").Select()
